$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Uren Gepland") and E ("Uren Besteed") were stored as text like
# "5 uur" / "1.5 uur" / "25 min" etc. They are replaced with plain numeric
# values representing the number of hours.

$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 5

$ws.Range("D4").Value = 8
$ws.Range("E4").Value = 3

$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

$ws.Range("D7").Value = 2

$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1

$ws.Range("D9").Value = 1.5
$ws.Range("E9").Value = 1

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

$ws.Range("D11").Value = 0.5
$ws.Range("E11").Value = 0.5

$ws.Range("D12").Value = 0.5
$ws.Range("E12").Value = 0.5

$ws.Range("D13").Value = 2

$ws.Range("D14").Value = 2

$ws.Range("D15").Value = 1

$ws.Range("D16").Value = 2

# Row 17: activity text changed, hours reworked, and the "Uren Besteed" cell
# is cleared entirely.
$ws.Range("B17").Value = "toevoegen items (admin)"
$ws.Range("D17").Value = 0.5
$ws.Range("E17").ClearContents()

$ws.Range("D18").Value = 0.5
$ws.Range("E18").Value = 0.5

# Update the active selection to match the saved view state.
$ws.Range("E19").Select()
